$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("G2").Value = 117.9639543333333
$ws.Range("H2").Value = 353.891863
$ws.Range("I2").Value = 0.2661690114309019
$ws.Range("J2").Value = 0.2661690114309019
$ws.Range("K2").Value = 3.0
$ws.Range("M2").Value = 2.618716333333334
$ws.Range("N2").Value = 7.856149000000001
$ws.Range("O2").Value = 0.07115908183301342
$ws.Range("P2").Value = 0.07115908183301341
$ws.Range("Q2").Value = 308.9141339572875
$ws.Range("R2").Value = 2780.227205615588
$ws.Range("S2").Value = 0.01894034246582384
$ws.Range("T2").Value = 0.01894034246582383
$ws.Range("E3").Value = 3.0
$ws.Range("G3").Value = 117.9639543333333
$ws.Range("H3").Value = 353.891863
$ws.Range("I3").Value = 0.2661690114309019
$ws.Range("J3").Value = 0.2661690114309019
$ws.Range("K3").Value = 3.0
$ws.Range("M3").Value = 15.503283
$ws.Range("N3").Value = 46.509849
$ws.Range("O3").Value = 0.4212748702999519
$ws.Range("P3").Value = 0.4212748702999519
$ws.Range("Q3").Value = 1828.828567828743
$ws.Range("R3").Value = 16459.45711045869
$ws.Range("S3").Value = 0.1121303157684196
$ws.Range("T3").Value = 0.1121303157684196
$ws.Range("E4").Value = 3.0
$ws.Range("G4").Value = 117.9639543333333
$ws.Range("H4").Value = 353.891863
$ws.Range("I4").Value = 0.2661690114309019
$ws.Range("J4").Value = 0.2661690114309019
$ws.Range("K4").Value = 3.0
$ws.Range("M4").Value = 18.67887366666666
$ws.Range("N4").Value = 56.036621
$ws.Range("O4").Value = 0.5075660478670347
$ws.Range("P4").Value = 0.5075660478670347
$ws.Range("Q4").Value = 2203.433800212769
$ws.Range("R4").Value = 19830.90420191492
$ws.Range("S4").Value = 0.1350983531966585
$ws.Range("T4").Value = 0.1350983531966585
$ws.Range("E5").Value = 3.0
$ws.Range("G5").Value = 282.6413673333333
$ws.Range("H5").Value = 847.924102
$ws.Range("I5").Value = 0.6377403483780447
$ws.Range("J5").Value = 0.6377403483780446
$ws.Range("K5").Value = 3.0
$ws.Range("M5").Value = 2.618716333333334
$ws.Range("N5").Value = 7.856149000000001
$ws.Range("O5").Value = 0.07115908183301342
$ws.Range("P5").Value = 0.07115908183301341
$ws.Range("Q5").Value = 740.1575651114666
$ws.Range("R5").Value = 6661.418086003198
$ws.Range("S5").Value = 0.04538101763844777
$ws.Range("T5").Value = 0.04538101763844776
$ws.Range("E6").Value = 3.0
$ws.Range("G6").Value = 282.6413673333333
$ws.Range("H6").Value = 847.924102
$ws.Range("I6").Value = 0.6377403483780447
$ws.Range("J6").Value = 0.6377403483780446
$ws.Range("K6").Value = 3.0
$ws.Range("M6").Value = 15.503283
$ws.Range("N6").Value = 46.509849
$ws.Range("O6").Value = 0.4212748702999519
$ws.Range("P6").Value = 0.4212748702999519
$ws.Range("Q6").Value = 4381.869105275622
$ws.Range("R6").Value = 39436.8219474806
$ws.Range("S6").Value = 0.2686639825480069
$ws.Range("T6").Value = 0.2686639825480069
$ws.Range("E7").Value = 3.0
$ws.Range("G7").Value = 282.6413673333333
$ws.Range("H7").Value = 847.924102
$ws.Range("I7").Value = 0.6377403483780447
$ws.Range("J7").Value = 0.6377403483780446
$ws.Range("K7").Value = 3.0
$ws.Range("M7").Value = 18.67887366666666
$ws.Range("N7").Value = 56.036621
$ws.Range("O7").Value = 0.5075660478670347
$ws.Range("P7").Value = 0.5075660478670347
$ws.Range("Q7").Value = 5279.42239339326
$ws.Range("R7").Value = 47514.80154053934
$ws.Range("S7").Value = 0.3236953481915901
$ws.Range("T7").Value = 0.32369534819159
$ws.Range("E8").Value = 3.0
$ws.Range("G8").Value = 42.586595
$ws.Range("H8").Value = 127.759785
$ws.Range("I8").Value = 0.09609064019105341
$ws.Range("J8").Value = 0.09609064019105343
$ws.Range("K8").Value = 3.0
$ws.Range("M8").Value = 2.618716333333334
$ws.Range("N8").Value = 7.856149000000001
$ws.Range("O8").Value = 0.07115908183301342
$ws.Range("P8").Value = 0.07115908183301341
$ws.Range("Q8").Value = 111.5222119075517
$ws.Range("R8").Value = 1003.699907167965
$ws.Range("S8").Value = 0.006837721728741818
$ws.Range("T8").Value = 0.006837721728741818
$ws.Range("E9").Value = 3.0
$ws.Range("G9").Value = 42.586595
$ws.Range("H9").Value = 127.759785
$ws.Range("I9").Value = 0.09609064019105341
$ws.Range("J9").Value = 0.09609064019105343
$ws.Range("K9").Value = 3.0
$ws.Range("M9").Value = 15.503283
$ws.Range("N9").Value = 46.509849
$ws.Range("O9").Value = 0.4212748702999519
$ws.Range("P9").Value = 0.4212748702999519
$ws.Range("Q9").Value = 660.232034291385
$ws.Range("R9").Value = 5942.088308622465
$ws.Range("S9").Value = 0.04048057198352537
$ws.Range("T9").Value = 0.04048057198352537
$ws.Range("E10").Value = 3.0
$ws.Range("G10").Value = 42.586595
$ws.Range("H10").Value = 127.759785
$ws.Range("I10").Value = 0.09609064019105341
$ws.Range("J10").Value = 0.09609064019105343
$ws.Range("K10").Value = 3.0
$ws.Range("M10").Value = 18.67887366666666
$ws.Range("N10").Value = 56.036621
$ws.Range("O10").Value = 0.5075660478670347
$ws.Range("P10").Value = 0.5075660478670347
$ws.Range("Q10").Value = 795.4696278984982
$ws.Range("R10").Value = 7159.226651086484
$ws.Range("S10").Value = 0.04877234647878623
$ws.Range("T10").Value = 0.04877234647878623
